$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: OriginalFilename header + 15 data rows
$values = @(
    "OriginalFilename",
    "1_TIG_ACnopulse_900_IQ_time.iq",
    "1_TIG_ACnopulse_900_IQ_time_2.iq",
    "1_TIG_ACnopulse_24_IQ_time.iq",
    "2_TIG_ACnopulse_900_IQ_time.iq",
    "2_TIG_ACnopulse_24GHz_IQ_time.iq",
    "3_TIG_ACpulse_900_IQ_time.iq",
    "3_TIG_ACpulse_24GHz_IQ_time.iq",
    "4_TIG_ACpulse_900_IQ_time.iq",
    "4_TIG_ACpulse_24GHz_IQ_time.iq",
    "5_TIG_ACnopulse_900_IQ_time.iq",
    "5_TIG_ACnopulse_24GHz_IQ_time.iq",
    "6_TIG_DCnopulse_900_IQ_time.iq",
    "6_TIG_DCnopulse_24GHz_IQ_time.iq",
    "11_TIG_DCnopulse_900_IQ_time.iq",
    "11_TIG_DCnopulse_24GHz_IQ_time.iq"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 11).Value = $values[$i]
}

# Set column K width to match the new layout (closest achievable value;
# the COM ColumnWidth setter quantizes to 1/6-character steps)
$ws.Range("K1:K16").ColumnWidth = 32.6
